$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting
# Late/heading(Office)/Outstanding columns one to the right.
$ws.Columns("N:N").Insert()

# Newly inserted column inherits the format of the column to its left (M);
# make the stored width match that as closely as this runtime allows.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Select the Repayment schedule sheet (making it the active tab) and move
# the selection to K15, matching the author's final cursor position.
$ws.Range("K15").Select()
